$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 48 - 四方坪站 station data for 2025-10-24 (date serial 45954)
$ws.Range("A48").Value = 45954
$ws.Range("B48").Value = "四方坪站"
$ws.Range("C48").Value = 9462.42
$ws.Range("D48").Value = 7788.69
$ws.Range("E48").Value = 3311.88
$ws.Range("F48").Value = 397

# New row 49 - 高岭站 station data for 2025-10-24 (date serial 45954)
$ws.Range("A49").Value = 45954
$ws.Range("B49").Value = "高岭站"
$ws.Range("C49").Value = 5405.33
$ws.Range("D49").Value = 4557.09
$ws.Range("E49").Value = 1398.83
$ws.Range("F49").Value = 200

# Move the selection to match the workbook's saved cursor position
[void]$ws.Range("G56").Select()
